# Apply the ELCData.xlsx edit: fill in cell A1 ("Cars") and add a new
# row of scraped data in row 7 (A7:D7), growing the sheet's used range
# from A1:B2 to A1:D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 was an empty text cell; it now holds the category label.
$ws.Range("A1").Value = "Cars"

# New row 7 with the scraped product-listing fields.
$ws.Range("A7").Value = "https://www.elc.co.uk/brands/paw-patrol"
$ws.Range("B7").Value = "Narrow Your Results"
$ws.Range("C7").Value = "Available to pick up for FREE from our stores within 30 minutes."
$ws.Range("D7").Value = "Add to Basket"
